# Insert a new week's worth of data (2 rows: "Primera" and "Segunda")
# right after the existing row 900, shifting all following rows down by 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("901:902").Insert()

# New row 901 - Primera
$ws.Cells.Item(901, 1).Value = 9
$ws.Cells.Item(901, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(901, 3).Value = "Metropolitana"
$ws.Cells.Item(901, 4).Value = 45147
$ws.Cells.Item(901, 5).Value = 13
$ws.Cells.Item(901, 6).Value = 100112009
$ws.Cells.Item(901, 7).Value = "Acelga"
$ws.Cells.Item(901, 8).Value = "Sin especificar"
$ws.Cells.Item(901, 9).Value = "Primera"
$ws.Cells.Item(901, 10).Value = 70
$ws.Cells.Item(901, 11).Value = 11000
$ws.Cells.Item(901, 12).Value = 11000
$ws.Cells.Item(901, 13).Value = 11000
$ws.Cells.Item(901, 14).Value = "`$/docena de atados"
$ws.Cells.Item(901, 15).Value = "Región Metropolitana"
$ws.Cells.Item(901, 16).Value = 3667
$ws.Cells.Item(901, 17).Value = 3
$ws.Cells.Item(901, 18).Value = "Hortaliza"

# New row 902 - Segunda
$ws.Cells.Item(902, 1).Value = 9
$ws.Cells.Item(902, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(902, 3).Value = "Metropolitana"
$ws.Cells.Item(902, 4).Value = 45147
$ws.Cells.Item(902, 5).Value = 13
$ws.Cells.Item(902, 6).Value = 100112009
$ws.Cells.Item(902, 7).Value = "Acelga"
$ws.Cells.Item(902, 8).Value = "Sin especificar"
$ws.Cells.Item(902, 9).Value = "Segunda"
$ws.Cells.Item(902, 10).Value = 52
$ws.Cells.Item(902, 11).Value = 9000
$ws.Cells.Item(902, 12).Value = 9000
$ws.Cells.Item(902, 13).Value = 9000
$ws.Cells.Item(902, 14).Value = "`$/docena de atados"
$ws.Cells.Item(902, 15).Value = "Región Metropolitana"
$ws.Cells.Item(902, 16).Value = 3000
$ws.Cells.Item(902, 17).Value = 3
$ws.Cells.Item(902, 18).Value = "Hortaliza"

Write-Host "Inserted new rows 901-902"
